$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared text blocks reused across the "Message" style sheets
# ---------------------------------------------------------------------------
$neo4jUrl = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userName = "neo4j"
$pwdValue = "icdcDBneo4j0"
$outputPath = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC11_Canine_Filter_Diagnosis-MaligLymph_Neo4jData.xlsx"

$cypherOutputQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN ['Malignant lymphoma'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

$statOutputQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Malignant lymphoma']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

function Set-MessageSheet($ws, $cypherText) {
    $ws.Range("A1").Value = "Neo4j_URL:"
    $ws.Range("A2").Value = $neo4jUrl
    $ws.Range("A3").Value = "User_name:"
    $ws.Range("A4").Value = $userName
    $ws.Range("A5").Value = "PWD:"
    $ws.Range("A6").Value = $pwdValue
    $ws.Range("A7").Value = "Cypher:"
    $ws.Range("A8").Value = $cypherText
    $ws.Range("A9").Value = "Output:"
    $ws.Range("A10").Value = $outputPath
}

# ---------------------------------------------------------------------------
# New sheet: CypherOutput_Message (repeats the same Message layout)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCypherMsg = $wb.Worksheets.Add($null, $lastSheet)
$wsCypherMsg.Name = "CypherOutput_Message"
Set-MessageSheet $wsCypherMsg $cypherOutputQuery

# ---------------------------------------------------------------------------
# New sheet: StatOutput (aggregated counts table)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStat = $wb.Worksheets.Add($null, $lastSheet)
$wsStat.Name = "StatOutput"
$wsStat.Range("A1").Value = "number_of_files"
$wsStat.Range("B1").Value = "number_of_sample"
$wsStat.Range("C1").Value = "number_of_cases"
$wsStat.Range("D1").Value = "number_of_study"
$wsStat.Range("A2").Value = "'0"
$wsStat.Range("B2").Value = "'0"
$wsStat.Range("C2").Value = "'6"
$wsStat.Range("D2").Value = "'1"

# ---------------------------------------------------------------------------
# New sheet: StatOutput_Message (Message layout repeated twice; the second
# block's Cypher text is the stat query)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStatMsg = $wb.Worksheets.Add($null, $lastSheet)
$wsStatMsg.Name = "StatOutput_Message"
Set-MessageSheet $wsStatMsg $cypherOutputQuery

$wsStatMsg.Range("A11").Value = "Neo4j_URL:"
$wsStatMsg.Range("A12").Value = $neo4jUrl
$wsStatMsg.Range("A13").Value = "User_name:"
$wsStatMsg.Range("A14").Value = $userName
$wsStatMsg.Range("A15").Value = "PWD:"
$wsStatMsg.Range("A16").Value = $pwdValue
$wsStatMsg.Range("A17").Value = "Cypher:"
$wsStatMsg.Range("A18").Value = $statOutputQuery
$wsStatMsg.Range("A19").Value = "Output:"
$wsStatMsg.Range("A20").Value = $outputPath
